$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert 3 new rows at row 22 (Bosses / Adding Content / NPC) ---
# Everything currently at rows 22:45 slides down to 25:48.
$ws.Rows("22:24").Insert()

# --- Fill in the 3 newly inserted rows ---
$ws.Range("A22").Value = "Bosses"
$ws.Range("B22").Value = "Not Started"
$ws.Range("C22").Value = "Alex"

$ws.Range("A23").Value = "Adding Content"
$ws.Range("B23").Value = "Started"
$ws.Range("C23").Value = "Alex, John, Matt"

$ws.Range("A24").Value = "NPC"
$ws.Range("B24").Value = "Started"
$ws.Range("C24").Value = "Matt"

# --- Newly populated "Who" cells for rows that already existed (now shifted) ---
$ws.Range("C25").Value = "Matt"          # Character Saves/Loads
$ws.Range("C26").Value = "John"          # Inventory

# --- "Who" annotations edited with '*' markers on existing Done/Started rows ---
$ws.Range("C6").Value  = "Alex*, John, Matt"   # Entities
$ws.Range("C7").Value  = "Alex*, John"         # Hud
$ws.Range("C11").Value = "Alex, John, Matt*"   # Hub World
$ws.Range("C12").Value = "John*, Matt"         # Combat
$ws.Range("C15").Value = "Alex, John*, Joe"    # Experience System
$ws.Range("C18").Value = "Alex, John, Matt*"   # Medieval World
$ws.Range("C19").Value = "Alex*, Matt"         # Dungeon

# --- Dynamic Monsters Lvls row gains a "Who" and a goal note in column D ---
$ws.Range("C34").Value = "Alex"
$ws.Range("D34").Value = "Goal***!*!*!*!*!*!*!*!*!**!!!"

# --- Networking rows gain "Who" = Alex ---
$ws.Range("C43").Value = "Alex"   # Basic Networking
$ws.Range("C44").Value = "Alex"   # Chat System (console)
$ws.Range("C45").Value = "Alex"   # Sync Ai

# --- Restore the current selection to match the saved view ---
$ws.Range("C24").Select()
